# إضافة حدث جديد في Card13
# 1) Every previously-blank data cell (rows 2-15) on the Card13 sheet gets the
#    literal text "nan" written into it (this mirrors the source workbook's
#    export/round-trip behaviour captured in the commit).
# 2) A brand new row (row 16) is appended with the new service event.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card13")

$nanCells = @(
    "D2","E2","F2","G2","H2","I2","J2","K2","L2","M2","N2","O2",
    "G3","H3","I3","J3","K3","M3","N3","O3",
    "D4","E4","F4","G4","H4","I4","J4","K4","L4","M4","N4","O4",
    "D5","E5","H5","I5","J5","K5","M5","N5","O5",
    "E6","F6","G6","I6","J6","K6","M6","N6","O6",
    "E7","G7","H7","I7","J7","M7","N7","O7",
    "D8","F8","G8","H8","J8","K8",
    "E9","H9","I9","J9","K9","M9","N9","O9",
    "D10","E10","F10","G10","H10","I10","J10","K10","L10","M10","N10","O10",
    "D11","E11","F11","G11","H11","I11","J11","K11","L11","M11","N11","O11",
    "D12","E12","F12","G12","H12","I12","J12","K12","L12","M12","N12","O12",
    "D13","E13","F13","G13","H13","I13","J13","K13","L13","M13","N13","O13",
    "B14","C14","D14","E14","F14","G14","H14","I14","J14","K14","N14",
    "B15","C15","D15","E15","F15","G15","H15","I15","J15","K15","N15"
)

foreach ($addr in $nanCells) {
    $ws.Range($addr).Value = "nan"
}

# New row 16 - the newly logged service event for Card13
$ws.Range("A16").Value = "13"
$ws.Range("L16").Value = "24\12\2024"
$ws.Range("M16").Value = "تم سن الفلاتس ومعايره"
$ws.Range("O16").Value = "الخبير"
